# Update cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell's value while forcing it to remain a text cell
# (prevents Excel from auto-converting numeric-looking strings like
# "311.73" into real numbers) and then restoring the cell's original
# "Normal" style so no stray formatting/style index is left behind.
function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

function Set-Cell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
Set-TextCell "D2" "27.458.84"
Set-TextCell "E2" "  +1.81%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.859.17"
Set-TextCell "E3" "  +0.83%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  -0.07%  "

# Row 5 - BNB
Set-TextCell "D5" "311.73"
Set-TextCell "E5" "  +0.94%  "

# Row 6 - USDC
Set-TextCell "E6" "  -0.06%  "

# Row 7 - XRP
Set-TextCell "D7" "0.4773"
Set-TextCell "E7" "  +0.10%  "

# Row 8 - Cardano
Set-TextCell "D8" "0.3799"
Set-TextCell "E8" "  +3.32%  "

# Row 9 - Dogecoin
Set-TextCell "D9" "0.07314"
Set-TextCell "E9" "  +1.55%  "

# Row 10 - Polygon
Set-TextCell "D10" "0.9297"
Set-TextCell "E10" "  +0.06%  "

# Row 11 - Solana
Set-TextCell "D11" "20.72"
Set-TextCell "E11" "  +4.66%  "

# Row 12 - TRON
Set-TextCell "D12" "0.07791"
Set-TextCell "E12" "  +0.90%  "

# Row 13 - WrappedEther
Set-TextCell "D13" "1.864.95"
Set-TextCell "E13" "  +0.75%  "

# Row 14 - Polkadot
Set-TextCell "D14" "5.443"
Set-TextCell "E14" "  +0.74%  "

# Row 15 - Chainlink
Set-TextCell "D15" "6.547"
Set-TextCell "E15" "  +1.61%  "

# Row 16 - Litecoin
Set-TextCell "D16" "90.22"
Set-TextCell "E16" "  +1.60%  "

# Row 17 - BinanceUSD
Set-TextCell "E17" "  -0.30%  "

# Row 18 - ShibaInu
Set-TextCell "D18" "0.000008817"

# Row 19 - Dai
Set-TextCell "E19" "  -0.16%  "

# Row 20 - WrappedBTC
Set-TextCell "D20" "27.452.32"
Set-TextCell "E20" "  +1.55%  "

# Row 21 - Avalanche
Set-TextCell "D21" "14.63"
Set-TextCell "E21" "  +0.60%  "

# Row 22 - Uniswap
Set-TextCell "D22" "5.096"
Set-TextCell "E22" "  +0.58%  "

# Row 23 - Cosmos
Set-TextCell "E23" "  +0.51%  "

# Row 24 - Toncoin
Set-TextCell "D24" "1.943"
Set-TextCell "E24" "  +0.02%  "

# Row 25 - Monero
Set-TextCell "D25" "154.86"
Set-TextCell "E25" "  +1.55%  "

# Row 26 - EthereumClassic
Set-TextCell "D26" "18.45"
Set-TextCell "E26" "  +1.44%  "

# Row 27 - LidoDAOToken
Set-TextCell "D27" "2.003"
Set-TextCell "E27" "  -0.48%  "

# Row 28 - BitcoinCash
Set-TextCell "D28" "115.39"
Set-TextCell "E28" "  +0.92%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextCell "D29" "4.940"

# Row 30 - Stellar
Set-TextCell "D30" "0.08895"
Set-TextCell "E30" "  +0.42%  "

# Row 31 - HuobiToken
Set-TextCell "E31" "  +0.47%  "

# Row 32 - ARBITRUM
Set-TextCell "D32" "1.205"
Set-TextCell "E32" "  +2.58%  "

# Row 33 - ImmutableX
Set-TextCell "D33" "0.7526"
Set-TextCell "E33" "  +1.79%  "

# Row 34 - Filecoin
Set-TextCell "D34" "4.579"
Set-TextCell "E34" "  +1.91%  "

# Row 35 - RenderToken
Set-TextCell "D35" "2.699"
Set-TextCell "E35" "  +0.25%  "

# Rows 36 & 37 swap: VeChain moves to row 36, TrustWalletToken moves to row 37
Set-Cell "B36" "VeChain"
Set-Cell "C36" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D36" "0.02048"
Set-TextCell "E36" "  +4.66%  "

Set-Cell "B37" "TrustWalletToken"
Set-Cell "C37" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D37" "1.124"
Set-TextCell "E37" "  +1.51%  "

# Row 38 - TheSandbox
Set-TextCell "D38" "0.5562"
Set-TextCell "E38" "  +6.03%  "

# Row 39 - Hedera
Set-TextCell "D39" "0.05275"
Set-TextCell "E39" "  +0.39%  "

# Row 40 - MXToken
Set-TextCell "D40" "2.988"
Set-TextCell "E40" "  +0.76%  "

# Row 41 - FraxShare
Set-TextCell "D41" "7.017"
Set-TextCell "E41" "  +0.18%  "

# Row 42 - Aptos
Set-TextCell "D42" "8.568"
Set-TextCell "E42" "  +3.45%  "

# Row 43 - Algorand
Set-TextCell "D43" "0.1515"
Set-TextCell "E43" "  +0.27%  "

# Row 44 - Decentraland
Set-TextCell "D44" "0.4869"
Set-TextCell "E44" "  +2.90%  "

# Row 45 - EnergySwap
Set-TextCell "E45" "  +0.87%  "

# Row 46 - PaxDollar
Set-TextCell "E46" "  -0.15%  "

# Row 47 - NEARProtocol
Set-TextCell "D47" "1.660"
Set-TextCell "E47" "  +3.66%  "

# Row 48 - Quant
Set-TextCell "D48" "103.12"
Set-TextCell "E48" "  +1.27%  "

# Row 49 - Aave
Set-TextCell "D49" "67.34"
Set-TextCell "E49" "  +2.39%  "

# Row 50 - Cronos
Set-TextCell "D50" "0.06102"
Set-TextCell "E50" "  +0.55%  "

# Row 51 - EOS
Set-TextCell "D51" "0.9137"
Set-TextCell "E51" "  +2.97%  "

Write-Host "Cryptos list updated"
